# edit.ps1 - apply the "Az Android" deck update:
#  - bump master title/body default font sizes (42->44, 20->28)
#  - add body bullet content to slide 2 ("Alapok") + mark its title run dirty
#  - append four new Title+Content slides (3-6): Tortenete, Kompatibilis
#    Eszkozok, Elonyok, Hatranyok

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Slide master default text sizes (titleStyle lvl1 42->44, bodyStyle
#    lvl1 20->28)
# ---------------------------------------------------------------------
$master = $p.SlideMaster
$titleStyle = $master.TextStyles.Item([Microsoft.Office.Interop.PowerPoint.PpTextStyleType]::ppTitleStyle)
$titleLvl1 = $titleStyle.Levels(1)
$titleLvl1.Font.Size = 44

$bodyStyle = $master.TextStyles.Item([Microsoft.Office.Interop.PowerPoint.PpTextStyleType]::ppBodyStyle)
$bodyLvl1 = $bodyStyle.Levels(1)
$bodyLvl1.Font.Size = 28

# ---------------------------------------------------------------------
# 2. Slide 2 ("Alapok") - fill in the empty content placeholder and
#    touch up the title run.
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

$titleShape = $s2.Shapes.Item(1)
# Re-assert the title text so the run carries dirty="0" like the target.
$titleShape.TextFrame.TextRange.Text = "Alapok"

$bodyShape = $s2.Shapes.Item(2)
$bodyTr = $bodyShape.TextFrame.TextRange

$bodyText = "Linux kernelt használ`rMobil eszközökön használható`rA fejlesztők Java nyelven írhatnak rá`rGoogle alkalmazások alapok`rRendszerrel összekapcsolt fiókok`rNyílt forrású OS`r"
$newRange = $bodyTr.InsertBefore($bodyText)
$newRange.LanguageID = "hu-HU"

# ---------------------------------------------------------------------
# 3. New slides 3-6 (Title and Content layout == index 2 on this master)
# ---------------------------------------------------------------------
function Set-TitleText($slide, [string]$text) {
    $sh = $slide.Shapes.Item(1)
    $sh.TextFrame.TextRange.Text = $text
    $sh.TextFrame.TextRange.LanguageID = "hu-HU"
}

function Set-BodyText($slide, [string]$text) {
    $sh = $slide.Shapes.Item(2)
    $tr = $sh.TextFrame.TextRange
    $tr.Text = $text
    $tr.LanguageID = "hu-HU"
}

# --- Slide 3: Tortenete -------------------------------------------------
$s3 = $p.Slides.Add(3, [Microsoft.Office.Interop.PowerPoint.PpSlideLayout]::ppLayoutText)
$s3.Shapes.Item(1).Name = "Cím 1"
$s3.Shapes.Item(2).Name = "Tartalom helye 2"
Set-TitleText $s3 "Története"
Set-BodyText $s3 "Az Android Inc. kezdte  fejleszteni`r2005-ben felvásárolta a Google`rKonkurens OS az iOS`rElső androidos telefon 2008-ban`rFolyamatos fejlesztés, legújabb a 12-es`r"

# --- Slide 4: Kompatibilis Eszkozok -------------------------------------
$s4 = $p.Slides.Add(4, [Microsoft.Office.Interop.PowerPoint.PpSlideLayout]::ppLayoutText)
$s4.Shapes.Item(1).Name = "Cím 1"
$s4.Shapes.Item(2).Name = "Tartalom helye 2"
Set-TitleText $s4 "Kompatibilis Eszközök"
Set-BodyText $s4 "Főleg érintőképernyős okoseszközök`rFőbb támogatók: Samsung, Xiaomi`r"

# --- Slide 5: Elonyok ----------------------------------------------------
$s5 = $p.Slides.Add(5, [Microsoft.Office.Interop.PowerPoint.PpSlideLayout]::ppLayoutText)
$s5.Shapes.Item(1).Name = "Cím 1"
$s5.Shapes.Item(2).Name = "Tartalom helye 2"
Set-TitleText $s5 "Előnyök"
Set-BodyText $s5 "Nyílt forráskód, több cég fejleszti`rBárhonnan letölthetünk fájlokat`rSD-kártya támogatott`rKönnyű használat"

# --- Slide 6: Hatranyok ---------------------------------------------------
$s6 = $p.Slides.Add(6, [Microsoft.Office.Interop.PowerPoint.PpSlideLayout]::ppLayoutText)
$s6.Shapes.Item(1).Name = "Cím 1"
$s6.Shapes.Item(2).Name = "Tartalom helye 2"
Set-TitleText $s6 "Hátrányok"
$s6Body = $s6.Shapes.Item(2)
$s6Tr = $s6Body.TextFrame.TextRange
$s6Tr.Text = "Nem elég tiszta a megjelenés (iOS-hez)`r"
$s6Tr.LanguageID = "hu-HU"
$s6Full = $s6Body.TextFrame.TextRange
$s6Last = $s6Full.Characters($s6Full.Length, 1)
$s6Last.IndentLevel = 2
